# Refactor handleSolve method to include GameMap parameter and improve puzzle
# solving logic -- corresponding bug-report workbook update: log a new bug
# entry (row 10) describing the Monolith/Strix puzzle spam issue and its fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new bug-report row (row 10)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Monolith/Strix puzzles can be spammed"
$ws.Range("C10").Value = "ricky"
$ws.Range("D10").Value = "no event triggering for puzzles"
$ws.Range("E10").Value = "After a correct answer, trigger event and reformat checking"
$ws.Range("F10").Value = "fixed"

# Widen column E to fit the new, longer description text
$ws.Columns.Item(5).ColumnWidth = 52.3

# Move the active selection to E7, as recorded when the workbook was saved
[void]$ws.Range("E7").Select()

# Page setup was touched (portrait orientation) when the sheet was saved
$ws.PageSetup.Orientation = 1
